$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '47.096.85'
$ws.Range('E2').Value = '  +3.43%  '
$ws.Range('D3').Value = '2.643.88'
$ws.Range('E3').Value = '  +9.83%  '
$ws.Range('E4').Value = '  -0.25%  '
$ws.Range('D5').Value = "'314.53"
$ws.Range('E5').Value = '  +5.02%  '
$ws.Range('D6').Value = "'103.65"
$ws.Range('E6').Value = '  +6.58%  '
$ws.Range('D7').Value = "'0.609"
$ws.Range('E7').Value = '  +8.07%  '
$ws.Range('D8').Value = "'0.999"
$ws.Range('E8').Value = '  -0.13%  '
$ws.Range('E9').Value = '  +15.71%  '
$ws.Range('D10').Value = "'39.31"
$ws.Range('E10').Value = '  +12.16%  '
$ws.Range('D11').Value = "'55.19"
$ws.Range('E11').Value = '  +2.52%  '
$ws.Range('D12').Value = "'0.0849"
$ws.Range('E12').Value = '  +6.96%  '
$ws.Range('D13').Value = "'8.37"
$ws.Range('E13').Value = '  +16.98%  '
$ws.Range('D14').Value = '3.039.24'
$ws.Range('E14').Value = '  +10.20%  '
$ws.Range('D16').Value = '2.652.41'
$ws.Range('E16').Value = '  +10.04%  '
$ws.Range('D17').Value = "'0.937"
$ws.Range('E17').Value = '  +10.83%  '
$ws.Range('D18').Value = "'15.21"
$ws.Range('E18').Value = '  +6.64%  '
$ws.Range('D19').Value = '47.493.54'
$ws.Range('E19').Value = '  +4.29%  '
$ws.Range('E20').Value = '  +8.55%  '
$ws.Range('D21').Value = "'13.37"
$ws.Range('E21').Value = '  +3.71%  '
$ws.Range('D22').Value = "'6.81"
$ws.Range('E22').Value = '  +9.24%  '
$ws.Range('D23').Value = "'71.65"
$ws.Range('E23').Value = '  +6.61%  '
$ws.Range('D24').Value = "'261.04"
$ws.Range('E24').Value = '  +7.72%  '
$ws.Range('D25').Value = "'3.13"
$ws.Range('E25').Value = '  +10.75%  '
$ws.Range('B26').Value = 'EthereumClassic'
$ws.Range('C26').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D26').Value = "'31.44"
$ws.Range('E26').Value = '  +47.93%  '
$ws.Range('B27').Value = 'ImmutableX'
$ws.Range('C27').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D27').Value = "'2.26"
$ws.Range('E27').Value = '  +17.27%  '
$ws.Range('E28').Value = '  +0.06%  '
$ws.Range('B29').Value = 'Cosmos'
$ws.Range('C29').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D29').Value = "'10.73"
$ws.Range('E29').Value = '  +9.71%  '
$ws.Range('B30').Value = 'InjectiveProtocol'
$ws.Range('C30').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D30').Value = "'41.40"
$ws.Range('E30').Value = '  +8.42%  '
$ws.Range('E31').Value = '  +4.16%  '
$ws.Range('D32').Value = "'6.31"
$ws.Range('E32').Value = '  +13.97%  '
$ws.Range('D33').Value = "'3.81"
$ws.Range('E33').Value = '  -1.39%  '
$ws.Range('D34').Value = "'2.33"
$ws.Range('E34').Value = '  +17.60%  '
$ws.Range('E35').Value = '  +5.22%  '
$ws.Range('D36').Value = "'0.0846"
$ws.Range('E36').Value = '  +8.97%  '
$ws.Range('D37').Value = "'153.33"
$ws.Range('E37').Value = '  +2.72%  '
$ws.Range('E38').Value = '  +4.38%  '
$ws.Range('E39').Value = '  +6.65%  '
$ws.Range('D40').Value = "'17.14"
$ws.Range('E40').Value = '  +11.76%  '
$ws.Range('D41').Value = "'4.37"
$ws.Range('E41').Value = '  +13.36%  '
$ws.Range('D42').Value = "'3.73"
$ws.Range('E43').Value = '  +10.13%  '
$ws.Range('D44').Value = "'21.70"
$ws.Range('E44').Value = '  +39.52%  '
$ws.Range('D45').Value = '2.062.20'
$ws.Range('E45').Value = '  +6.42%  '
$ws.Range('D46').Value = "'0.999"
$ws.Range('E46').Value = '  -0.06%  '
$ws.Range('D47').Value = "'93.85"
$ws.Range('E47').Value = '  +2.77%  '
$ws.Range('D48').Value = "'114.74"
$ws.Range('E48').Value = '  +11.78%  '
$ws.Range('E49').Value = '  +6.05%  '
$ws.Range('D50').Value = "'9.31"
$ws.Range('E50').Value = '  +6.52%  '
$ws.Range('D51').Value = "'0.203"
$ws.Range('E51').Value = '  +7.51%  '

# Cells whose values look like plain numbers (e.g. "314.53") were written
# with a leading quote to force text storage (matching the source sheet,
# which stores these as text). Strip the resulting quote-prefix styling so
# the cells keep the workbook default style, same as before the edit.
$numericTextCells = @('D5','D6','D7','D8','D10','D11','D12','D13','D17','D18','D21','D22','D23','D24','D25','D26','D27','D29','D30','D32','D33','D34','D36','D37','D40','D41','D42','D44','D46','D47','D48','D50','D51')
foreach ($addr in $numericTextCells) {
    $ws.Range($addr).ClearFormats()
}
